$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# The table "Tabla1" had columns: randomQTable | Epocas Training | Politica q | alpha | gamma
# It is being reshuffled to:      RandomQTable | Epocas Training | alpha | gamma | Politica
# (alpha/gamma shift left into the old "Politica q" slot, and a
# renamed "Politica " column is appended as the new last column.)
# -----------------------------------------------------------------

# Capture the original column C ("Politica q") values before we overwrite anything.
$politicaValues = @()
for ($r = 2; $r -le 6; $r++) {
    $politicaValues += $ws.Cells.Item($r, 3).Value2
}

# Capture alpha (D) / gamma (E) values too.
$alphaValues = @()
$gammaValues = @()
for ($r = 2; $r -le 6; $r++) {
    $alphaValues += $ws.Cells.Item($r, 4).Value2
    $gammaValues += $ws.Cells.Item($r, 5).Value2
}

# Rewrite header row: A=RandomQTable, B=Epocas Training, C=alpha, D=gamma, E=Politica (trailing space)
$ws.Range("A1").Value2 = "RandomQTable"
$ws.Range("C1").Value2 = "alpha"
$ws.Range("D1").Value2 = "gamma"
$ws.Range("E1").Value2 = "Política "

# Rewrite the data rows: alpha -> C, gamma -> D, Politica -> E
for ($i = 0; $i -lt 5; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 3).Value2 = $alphaValues[$i]
    $ws.Cells.Item($r, 4).Value2 = $gammaValues[$i]
    $ws.Cells.Item($r, 5).Value2 = $politicaValues[$i]
}

# Give the relocated "Politica " header the same look-and-feel as the other
# header cells (bold / bordered / centered "Encabezado 1" cell style).
$ws.Range("E1").Style = "Encabezado 1"

# Leftover formatted-but-empty column to the right of the table (F), matching
# the stray formatting left behind by the column shuffle.
$ws.Range("F1").Style = "Encabezado 1"
$ws.Range("F2:F6").Style = $ws.Range("A2").Style

# Column widths follow the data: alpha's old width now belongs to column C,
# gamma's old width to column D, and Politica's old width to column E.
$ws.Columns.Item(2).ColumnWidth = 24.5
$ws.Columns.Item(3).ColumnWidth = 23
$ws.Columns.Item(4).ColumnWidth = 28.16666666667
$ws.Columns.Item(5).ColumnWidth = 22.33333333333

# Final cursor position, mirroring the author's last selection.
$ws.Range("G8").Select()

Write-Output "done"
